$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 83, pushing existing rows 83-95 down to 84-96.
$ws.Rows.Item(83).Insert()

# Populate the new row 83 with the latest weekly record.
$ws.Cells.Item(83, 1).Value = 2
$ws.Cells.Item(83, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(83, 3).Value = "Coquimbo"
$ws.Cells.Item(83, 4).NumberFormat = $ws.Cells.Item(84, 4).NumberFormat
$ws.Cells.Item(83, 4).Value = 45127
$ws.Cells.Item(83, 5).Value = 4
$ws.Cells.Item(83, 6).Value = 100112022
$ws.Cells.Item(83, 7).Value = "Arveja Verde"
$ws.Cells.Item(83, 8).Value = "Perfection"
$ws.Cells.Item(83, 9).Value = "Primera"
$ws.Cells.Item(83, 10).Value = 700
$ws.Cells.Item(83, 11).Value = 20000
$ws.Cells.Item(83, 12).Value = 22000
$ws.Cells.Item(83, 13).Value = 21000
$ws.Cells.Item(83, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(83, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(83, 16).Value = 840
$ws.Cells.Item(83, 17).Value = 25
$ws.Cells.Item(83, 18).Value = "Hortaliza"
